$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the event name in row 4 from "Shadow Assembly 6PM" to "Shadow Assembly 7PM"
$ws.Range("A4").Value = "Shadow Assembly 7PM"

# Move the active selection to A5 (matches the saved selection state in the file)
$ws.Range("A5").Select()
